# Apply the diff: insert a new row at row 87 (pushing existing rows 87-156
# down to 88-157) and populate the new row 87 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 87; this shifts rows 87:156 down
# to 88:157 automatically (values, formats, styles all move together).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record values.
$ws.Range("A87").Value = 5
$ws.Range("B87").Value = "Macroferia Regional de Talca"
$ws.Range("C87").Value = "Maule"
$ws.Range("D87").Value = 44447
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112003
$ws.Range("G87").Value = "Ajo"
$ws.Range("H87").Value = "Chino"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 200
$ws.Range("K87").Value = 14000
$ws.Range("L87").Value = 14000
$ws.Range("M87").Value = 14000
$ws.Range("N87").Value = "$/caja 10 kilos"
$ws.Range("O87").Value = "China"
$ws.Range("P87").Value = 1400
$ws.Range("Q87").Value = 10
$ws.Range("R87").Value = "Hortaliza"
